$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 746, pushing the existing rows 746-787 down to 747-788
# (dimension grows from A1:D787 to A1:D788).
$ws.Rows.Item(746).EntireRow.Insert()

# The new row records another "2026/01/31" (Saturday) entry, at time 16, ranking 201.
# Column A/B already hold that exact date/weekday as plain text one row above (row 745),
# so clone them with Copy instead of re-typing the literal "2026/01/31" string - a plain
# Value assignment would get auto-parsed into a date serial number by the date-recognizer,
# which would not match how the rest of the date column is stored (as literal text).
$ws.Range("A745:B745").Copy($ws.Range("A746:B746"))
$ws.Range("C746").Value = 16
$ws.Range("D746").Value = 201
